# Fine tune the Argument.
# Paragraph-level rewrites: rephrase four body paragraphs of the essay
# to read more naturally (contractions expanded, sentences reworded,
# a few words swapped for synonyms), per the commit's intent.

$d = $word.ActiveDocument

# Paragraph 1 ("Ever since I was a teenager...")
$found = $d.Content.Find.Execute("Ever since I was a teenager, I’ve caught myself daydreaming about the future—wondering what life might look like, what crazy inventions might catch on, and which ones we’ll laugh about later. Out of everything, artificial intelligence has probably been the one thing I kept coming back to. It’s not just something out of sci-fi anymore, it’s in our phones, our homes, and even quietly influencing decisions we don’t realize we’re making.", $true, $false, $false, $false, $false, $true, 1, $false, "Ever since I was a teenager, I’ve caught myself daydreaming about the future — wondering what life might look like, what crazy inventions might catch on, and which ones we will laugh about later. Out of everything, artificial intelligence has probably been the one thing I kept coming back to. It is not just something out of sci-fi anymore, it is in our phones, our homes, and even quietly influencing decisions we do not even realize we are making.", 2)
if (-not $found) { throw "Replacement for p1 failed to match" }

# Paragraph 2 ("I believe my interest in AI began...")
$found = $d.Content.Find.Execute("I believe my interest in AI began after I read an article about self-driving cars. The idea that a machine could make split-second decisions on the road without a human behind the wheel honestly blew my mind. And yeah, it also freaked me out a bit. After that, I started seeing AI everywhere — in my phone’s recommendations, in ads that weirdly seemed to know what I was thinking about, even in apps I used for school projects. That’s when it hit me: AI isn’t just about robots or futuristic tech. It’s already woven into our daily lives.", $true, $false, $false, $false, $false, $true, 1, $false, "I believe my interest in AI began after I read an article about self-driving cars. The idea that a machine could make split-second decisions on the road without a human behind the wheel honestly blew my mind because of how shocking this whole concept is. After that, I started noticing AI everywhere — in phones’ recommendations, in ads that weirdly seemed to know exactly what we were thinking about. That is when it hit me: AI is not just about robots or futuristic tech. It’s already woven into our daily lives.", 2)
if (-not $found) { throw "Replacement for p2 failed to match" }

# Paragraph 3 ("One of the reasons I picked this topic...")
$found = $d.Content.Find.Execute("One of the reasons I picked this topic for my project is that I genuinely think AI is going to keep changing the way we live, in ways we’re only starting to notice. I’ve seen it already—some of the jobs in the market have started shifting because certain tasks are now automated. Meanwhile, I’m over here casually talking to voice assistants like they are roommates. Even the music I listen to and the games I play are shaped by algorithms that somehow know what I like before I do. What really fascinates me, though, isn’t just the tech itself, it’s how people respond to it. Some folks are all in, while others are more hesitant or even suspicious.", $true, $false, $false, $false, $false, $true, 1, $false, "One of the reasons I picked this topic for my project is that I genuinely think AI is going to keep changing the way we live, in ways we are only starting to notice. Some of them can be seen already — certain jobs in the market start to change because tasks can now be automated. Even the music we listen to and the games we all play are shaped by algorithms that somehow know what we like before we even do. What really fascinates me, though, is not just the tech itself, it’s how people respond to it. Some people are all in for it, while others are more hesitant or even skeptical.", 2)
if (-not $found) { throw "Replacement for p3 failed to match" }

# Paragraph 4 ("With this project, I want to explore...")
$found = $d.Content.Find.Execute("With this project, I want to explore how AI is affecting everyday life — not just the big headlines, but the small stuff too. From what shows up in our feeds to how we think about work and relationships, AI is already part of the picture. I’m hoping to capture both the excitement and the weirdness that comes with it. Because, like every major shift in technology, how we deal with AI now is going to shape our future in ways we probably can’t fully imagine yet.", $true, $false, $false, $false, $false, $true, 1, $false, "With this project, I want to explore how AI is affecting everyday life — not just the big headlines, but the small stuff too. From the content we see in our feeds to the way we approach work and relationships, AI is already part of the picture. I’m hoping to capture both the excitement and the weirdness that comes with it. Because, like every major shift in technology, how we deal with AI now is going to shape our future in ways we probably can’t fully comprehend yet.", 2)
if (-not $found) { throw "Replacement for p4 failed to match" }

